$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns B and C (values tuned so the engine's pixel-snap rounding
# lands as close as possible to the authored OOXML widths of 26.875 / 135.5)
$ws.Columns.Item(2).ColumnWidth = 26.15
$ws.Columns.Item(3).ColumnWidth = 134.85

# New row of data: "Bloodstained Chivalry"
$ws.Range("B4").Value = "Bloodstained Chivalry"
$ws.Range("C4").Value = ",`n[StatisticBoost;DamageModifier_PhysicalDamage=0.25;],`n,`n[StatisticBoost,Triggerable;DamageModifier_ChargedAttack=0.25;],"

# Match the formatting used by row 3 (vertical-centered; C wraps text)
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("C4").WrapText = $true

$ws.Rows.Item(4).RowHeight = 57

# Selection moves as part of the edit session
[void]$ws.Range("C22").Select()
